$d = $word.ActiveDocument

# Replace the placeholder "{name}" with "${name}" (prefix it with a literal $)
$d.Content.Find.Execute("{name}", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "`${name}", 2)
